# --------------------------------------------------------------------
# B1--and-B2-PowerPoint.pptx edit
#
# 1. The table on slide 5 (the "types of financial documents" table)
#    switches its built-in table style from
#       {F98CD05B-CCEC-4A68-A7DC-2B4E70FDE13E}
#    to
#       {55225FC8-5A2C-4A22-8BC5-6B7BF667BFF2}
#
# 2. The deck's theme colours (the ones actually driving every slide,
#    i.e. the slide-master theme) swap from the "Integral / Red Violet"
#    palette to the stock "Office Theme" palette (dk1/lt1/dk2/lt2 +
#    accent1-6 + hyperlink/followed-hyperlink).
# --------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 ---------------------------------------
$tableSlide  = $p.Slides.Item(5)
$tableShape  = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{55225FC8-5A2C-4A22-8BC5-6B7BF667BFF2}", $false)

# --- 2. Theme colours --------------------------------------------------
# The 12-slot DrawingML theme colour scheme, reached off any slide,
# is shared by the whole deck (it is backed by the slide master's
# theme part) so a single edit re-colours every slide consistently.
$tcs = $p.Slides.Item(1).ThemeColorScheme

$tcs.Colors(1).RGB  = 0         # dk1      000000
$tcs.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388   # dk2      44546A
$tcs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407     # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308  # accent5  4472C4
$tcs.Colors(10).RGB = 4697456   # accent6  70AD47
$tcs.Colors(11).RGB = 12673797  # hlink    0563C1
$tcs.Colors(12).RGB = 7491477   # folHlink 954F72
